$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, id (col A), topic (col B), seq code (col C, always text),
# subTopic (col D), subTopic-quotePrefix-flag, media URL (col E)
$rows = @(
    @(85, 84, "JavaScript - Variables, String Literals and String Concatenation", "1001", "Create file for learning Variables", 0, "https://www.youtube.com/embed/g9UWXxRRgnI"),
    @(86, 85, "JavaScript - Variables, String Literals and String Concatenation", "1002", "Define Variable", 0, "https://www.youtube.com/embed/o8EzyBoy3ZU"),
    @(87, 86, "JavaScript - Variables, String Literals and String Concatenation", "1003", "Use defined Variable", 0, "https://www.youtube.com/embed/pebkQDQ9MjQ"),
    @(88, 87, "JavaScript - Variables, String Literals and String Concatenation", "1004", "Change variable value", 0, "https://www.youtube.com/embed/kswH7btOgKs"),
    @(89, 88, "JavaScript - Variables, String Literals and String Concatenation", "1005", "String Concatenation", 0, "https://www.youtube.com/embed/e4a-cMZGjGo"),
    @(90, 89, "JavaScript - DOM and String Concatenation", "1101", "Create file for menu names display", 0, "https://www.youtube.com/embed/DskDYjNAOQU"),
    @(91, 90, "JavaScript - DOM and String Concatenation", "1102", "The 'let' variable should be unique", 0, "https://www.youtube.com/embed/HyQwkexwnE0"),
    @(92, 91, "JavaScript - DOM and String Concatenation", "1103", "Define unique variables", 0, "https://www.youtube.com/embed/aH8f-LsZKjM"),
    @(93, 92, "JavaScript - DOM and String Concatenation", "1104", "Concatenate with div tags", 0, "https://www.youtube.com/embed/5OwqMcYdjxM"),
    @(94, 93, "JavaScript - DOM and String Concatenation", "1105", "Concatenate with div tags in new variable", 0, "https://www.youtube.com/embed/bh65GkT3amI"),
    @(95, 94, "JavaScript - DOM and String Concatenation", "1106", "Concatenation explained", 0, "https://www.youtube.com/embed/fcmjbaHd1uY"),
    @(96, 95, "JavaScript - DOM and String Concatenation", "1107", "Include third menu item", 0, "https://www.youtube.com/embed/ksIvGO5oMfg"),
    @(97, 96, "JavaScript - Template Literal", "1201", "Create file for Template Literal", 0, "https://www.youtube.com/embed/-aQNxKojEe0"),
    @(98, 97, "JavaScript - Template Literal", "1202", "Block comment concatenation code", 0, "https://www.youtube.com/embed/BVCQ8Of62bc"),
    @(99, 98, "JavaScript - Template Literal", "1203", "Define empty string with backtick", 0, "https://www.youtube.com/embed/7aU8RTFSPAI"),
    @(100, 99, "JavaScript - Template Literal", "1204", "Define div tags within backtick", 0, "https://www.youtube.com/embed/rS-zMUVd2yU"),
    @(101, 100, "JavaScript - Template Literal", "1205", "Include template placeholders", 0, "https://www.youtube.com/embed/cX-drhxpaBM"),
    @(102, 101, "JavaScript - Template Literal", "1206", "Update menu names page with template literals", 0, "https://www.youtube.com/embed/nWXIGU-SndA"),
    @(103, 102, "JavaScript - Array", "1301", "Create file for array demo", 0, "https://www.youtube.com/embed/c0B2_R4AS4U"),
    @(104, 103, "JavaScript - Array", "1302", "Declare array of numbers", 0, "https://www.youtube.com/embed/4jk-sftdJIs"),
    @(105, 104, "JavaScript - Array", "1303", "Read array item", 0, "https://www.youtube.com/embed/rjBpMQDmi8c"),
    @(106, 105, "JavaScript - Array", "1304", "Array length", 0, "https://www.youtube.com/embed/CqtWdqg72OM"),
    @(107, 106, "JavaScript - 'for' loop", "1401", "Create file to learn 'for' loop", 0, "https://www.youtube.com/embed/Xbkg9vjcSQw"),
    @(108, 107, "JavaScript - 'for' loop", "1402", "Print numbers from 1 to 5 using for loop", 0, "https://www.youtube.com/embed/C5ksx7ngvno"),
    @(109, 108, "JavaScript - 'for' loop", "1403", "The 'for' loop construct", 1, "https://www.youtube.com/embed/WoKxUbqSoUg"),
    @(110, 109, "JavaScript - 'for' loop", "1404", "Initialization section of 'for' loop", 0, "https://www.youtube.com/embed/jaZpnAh07oc"),
    @(111, 110, "JavaScript - 'for' loop", "1405", "End statement of 'for' loop", 0, "https://www.youtube.com/embed/5wgut5SP2aw"),
    @(112, 111, "JavaScript - 'for' loop", "1406", "End condition of 'for' loop", 0, "https://www.youtube.com/embed/A0oubuO9_ZY"),
    @(113, 112, "JavaScript - 'for' loop", "1407", "Code execution flow of 'for' loop", 0, "https://www.youtube.com/embed/RrUJNeV7cR0"),
    @(114, 113, "JavaScript - 'for' loop", "1408", "Print array index numbers", 0, "https://www.youtube.com/embed/FldJzVXVRgE"),
    @(115, 114, "JavaScript - 'for' loop", "1409", "Print array items using 'for' loop", 0, "https://www.youtube.com/embed/ZPqjUuDu1pA"),
    @(116, 115, "JavaScript - 'for' loop", "1410", "Implement 'for .. of' loop", 0, "https://www.youtube.com/embed/vZ4DDzDE8LQ"),
    @(117, 116, "JavaScript - 'for' loop", "1411", "Display menu names using array and for", 0, "https://www.youtube.com/embed/UZfLeKz9Fhk"),
    @(118, 117, "JavaScript - 'for' loop", "1412", "Implement menu names as array declaration", 0, "https://www.youtube.com/embed/-SBDuMCu55s"),
    @(119, 118, "JavaScript - 'for' loop", "1413", "Implement menu names using 'for of' and Template Literals", 0, "https://www.youtube.com/embed/rZGpuaDw5WQ")

)

foreach ($entry in $rows) {
    $r       = $entry[0]
    $idVal   = $entry[1]
    $topic   = $entry[2]
    $seq     = $entry[3]
    $subTopic = $entry[4]
    $subTopicQuote = $entry[5]
    $media   = $entry[6]

    $ws.Cells.Item($r, 1).Value = $idVal
    $ws.Cells.Item($r, 2).Value = $topic
    $ws.Cells.Item($r, 3).Value = "'" + $seq
    if ($subTopicQuote -eq 1) {
        $ws.Cells.Item($r, 4).Value = "'" + $subTopic
    } else {
        $ws.Cells.Item($r, 4).Value = $subTopic
    }
    $ws.Cells.Item($r, 5).Value = $media
}

$ws.Range("E92").Select() | Out-Null
